# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across
# ALC, ARM, BSM, CRP, CUL, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 11993.333
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 11993.333
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H116").Value = 13673
$ws.Range("I116").Value = 36583.332
$ws.Range("J116").Value = 6799.9
$ws.Range("K116").Value = 36583.332
$ws.Range("L116").Value = 6799.9
$ws.Range("M116").Value = -33141.332
$ws.Range("N116").Value = -13683.9

$ws.Range("H127").Value = 4601.25
$ws.Range("I127").Value = 3385
$ws.Range("J127").Value = 8250
$ws.Range("K127").Value = 10155
$ws.Range("L127").Value = 24750
$ws.Range("M127").Value = -5195
$ws.Range("N127").Value = -34670

$ws.Range("H131").Value = 3345.5
$ws.Range("J131").Value = 4688.6665
$ws.Range("L131").Value = 14065.9995
$ws.Range("N131").Value = -24145.9995

$ws.Range("H132").Value = 1240.7059
$ws.Range("I132").Value = 1262
$ws.Range("K132").Value = 3786
$ws.Range("M132").Value = -1256

$ws.Range("H137").Value = 1776.8462
$ws.Range("I137").Value = 1554.4546
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 4663.3638
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -2113.3638
$ws.Range("N137").Value = -14100

$ws.Range("H138").Value = 2572.9268
$ws.Range("I138").Value = 2514.3635
$ws.Range("J138").Value = 2640.7368
$ws.Range("K138").Value = 7543.0905
$ws.Range("L138").Value = 7922.2104
$ws.Range("M138").Value = -2403.0905
$ws.Range("N138").Value = -18202.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4358.6123
$ws.Range("I32").Value = 2388.2104
$ws.Range("K32").Value = 2388.2104
$ws.Range("M32").Value = -2101.2104

$ws.Range("H61").Value = 4191.1577
$ws.Range("I61").Value = 2808.875
$ws.Range("K61").Value = 2808.875
$ws.Range("M61").Value = -2596.875

$ws.Range("H63").Value = 9001.666999999999
$ws.Range("I63").Value = 9001.666999999999
$ws.Range("K63").Value = 9001.666999999999
$ws.Range("M63").Value = -8315.666999999999

$ws.Range("H66").Value = 9001.666999999999
$ws.Range("I66").Value = 9001.666999999999
$ws.Range("K66").Value = 45008.335
$ws.Range("M66").Value = -41576.335

$ws.Range("H74").Value = 830.6774
$ws.Range("I74").Value = 791.4138
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 791.4138
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = 82.58619999999996
$ws.Range("N74").Value = -3148

$ws.Range("H77").Value = 830.6774
$ws.Range("I77").Value = 791.4138
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 3957.069
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = 410.9309999999996
$ws.Range("N77").Value = -15736

$ws.Range("H136").Value = 4191.1577
$ws.Range("I136").Value = 2808.875
$ws.Range("K136").Value = 8426.625
$ws.Range("M136").Value = -5876.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H105").Value = 2470
$ws.Range("I105").Value = 2470
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2470
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -723
$ws.Range("N105").ClearContents()

$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws.Range("H134").Value = 6931.61
$ws.Range("I134").Value = 6499.7837
$ws.Range("J134").Value = 10926
$ws.Range("K134").Value = 19499.3511
$ws.Range("L134").Value = 32778
$ws.Range("M134").Value = -16964.3511
$ws.Range("N134").Value = -37848

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8100
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -326

$ws.Range("H31").Value = 5035.3335
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5035.3335
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5035.3335
$ws.Range("N31").Value = -5625.3335
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 5035.3335
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5035.3335
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5035.3335
$ws.Range("N34").Value = -5439.3335
$ws.Range("M34").ClearContents()

$ws.Range("H132").Value = 2744.2222
$ws.Range("I132").Value = 1203.2222
$ws.Range("K132").Value = 3609.6666
$ws.Range("M132").Value = -1079.6666

$ws.Range("H134").Value = 954.7857
$ws.Range("I134").Value = 950.2308
$ws.Range("J134").Value = 1014
$ws.Range("K134").Value = 2850.6924
$ws.Range("L134").Value = 3042
$ws.Range("M134").Value = -315.6923999999999
$ws.Range("N134").Value = -8112

$ws.Range("H135").Value = 32308.889
$ws.Range("J135").Value = 32308.889
$ws.Range("L135").Value = 32308.889
$ws.Range("N135").Value = -42448.889

$ws.Range("H138").Value = 107000
$ws.Range("J138").Value = 107000
$ws.Range("L138").Value = 107000
$ws.Range("N138").Value = -117280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1430.375
$ws.Range("J107").Value = 1430.375
$ws.Range("L107").Value = 4291.125
$ws.Range("N107").Value = -8131.125

$ws.Range("H113").Value = 6601.8237
$ws.Range("I113").Value = 25351.75
$ws.Range("J113").Value = 832.61536
$ws.Range("K113").Value = 76055.25
$ws.Range("L113").Value = 2497.84608
$ws.Range("M113").Value = -73885.25
$ws.Range("N113").Value = -6837.84608

$ws.Range("H131").Value = 811.37
$ws.Range("I131").Value = 468.33334
$ws.Range("J131").Value = 833.2659
$ws.Range("K131").Value = 1405.00002
$ws.Range("L131").Value = 2499.7977
$ws.Range("M131").Value = 3634.99998
$ws.Range("N131").Value = -12579.7977

$ws.Range("H132").Value = 1348
$ws.Range("J132").Value = 1196
$ws.Range("L132").Value = 10764
$ws.Range("N132").Value = -15824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3420.2
$ws.Range("I22").Value = 5250.5
$ws.Range("K22").Value = 5250.5
$ws.Range("M22").Value = -4955.5

$ws.Range("H27").Value = 3420.2
$ws.Range("I27").Value = 5250.5
$ws.Range("K27").Value = 5250.5
$ws.Range("M27").Value = -5143.5

$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15820

$ws.Range("H40").Value = 8498
$ws.Range("I40").Value = 1996
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 1996
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -1860
$ws.Range("N40").Value = -15272

$ws.Range("H61").Value = 3750

$ws.Range("H113").Value = 3750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 12028400
$ws.Range("I11").Value = 30000000
$ws.Range("K11").Value = 30000000
$ws.Range("M11").Value = -29999858

$ws.Range("H17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("L17").Value = 1000
$ws.Range("N17").Value = -1344

$ws.Range("H39").Value = 19900
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H126").Value = 7660.6523
$ws.Range("I126").Value = 7096.5
$ws.Range("K126").Value = 21289.5
$ws.Range("M126").Value = -18819.5

$ws.Range("H132").Value = 1631.8846
$ws.Range("I132").Value = 1315.0454
$ws.Range("J132").Value = 3374.5
$ws.Range("K132").Value = 3945.1362
$ws.Range("L132").Value = 10123.5
$ws.Range("M132").Value = -1415.1362
$ws.Range("N132").Value = -15183.5
